# Update "paises.xlsx" (Pais sheet) with the latest COVID country stats
# and re-sort the two tied pairs (Eslovaquia/Mozambique, Montserrat/Islas
# Malvinas) plus refresh the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 10:39"

# --- Rusia (row 7) -----------------------------------------------------------
$ws.Range("B7").Value = 942106
$ws.Range("C7").Value = 4785
$ws.Range("D7").Value = 755513
$ws.Range("E7").Value = 170494
$ws.Range("G7").Value = 110
$ws.Range("H7").Value = 16099

# --- Filipinas (row 25) ------------------------------------------------------
$ws.Range("B25").Value = 178022
$ws.Range("C25").Value = 4339
$ws.Range("D25").Value = 114114
$ws.Range("E25").Value = 61025
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = 2883

# --- Polonia (row 47) ---------------------------------------------------------
$ws.Range("D47").Value = 40481
$ws.Range("E47").Value = 16217

# --- Singapur (row 49) --------------------------------------------------------
$ws.Range("B49").Value = 56099
$ws.Range("C49").Value = 68
$ws.Range("E49").Value = 3262

# --- Afganistan (row 60) -------------------------------------------------------
$ws.Range("B60").Value = 37856
$ws.Range("C60").Value = 257
$ws.Range("D60").Value = 27681
$ws.Range("E60").Value = 8790
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 1385

# --- Eslovaquia / Mozambique swap (rows 123-124) ------------------------------
# Eslovaquia's total now overtakes Mozambique's (unchanged) total, so it moves
# up to row 123 while Mozambique drops to row 124.
$ws.Range("A123").Value = "Eslovaquia"
$ws.Range("B123").Value = 3102
$ws.Range("C123").Value = 80
$ws.Range("D123").Value = 2014
$ws.Range("E123").Value = 1055
$ws.Range("H123").Value = 33

$ws.Range("A124").Value = "Mozambique"
$ws.Range("B124").Value = 3045
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 1291
$ws.Range("E124").Value = 1735
$ws.Range("H124").Value = 19

# --- Lituania (row 128) --------------------------------------------------------
$ws.Range("B128").Value = 2528
$ws.Range("C128").Value = 32
$ws.Range("D128").Value = 1747
$ws.Range("E128").Value = 699
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 82

# --- Estonia (row 133) ---------------------------------------------------------
$ws.Range("B133").Value = 2227
$ws.Range("C133").Value = 20
$ws.Range("D133").Value = 2009
$ws.Range("E133").Value = 155

# --- Montserrat / Islas Malvinas swap (rows 213-214) --------------------------
# Both are tied on total cases (13); Montserrat's recovered/deaths now sort it
# ahead of Islas Malvinas.
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("H214").Value = 0
